$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0004821219520434317
$ws.Range("E2").Value = 0.0004821219520434317

# Row 3
$ws.Range("D3").Value = 0.8843835049236194
$ws.Range("E3").Value = 0.8843835049236194

# Row 4
$ws.Range("D4").Value = 0.0001426636193318267
$ws.Range("E4").Value = 0.0001426636193318267

# Row 5
$ws.Range("D5").Value = 0.000000000000000008944065204798246
$ws.Range("E5").Value = 0.000000000000000008944065204798246

# Row 6
$ws.Range("D6").Value = 0.5723590334069846
$ws.Range("E6").Value = 0.5723590334069846

# Row 7
$ws.Range("D7").Value = 0.9999999770188902
$ws.Range("E7").Value = 0.00000002298110979293

# Row 8
$ws.Range("D8").Value = 0.9903156112971674
$ws.Range("E8").Value = 0.009684388702832614

# Row 9
$ws.Range("D9").Value = 0.9999980311993291
$ws.Range("E9").Value = 0.000001968800670937298

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.2812967851856813
$ws.Range("E10").Value = 0.7187032148143186

# Row 11
$ws.Range("D11").Value = 0.9556071031131752
$ws.Range("E11").Value = 0.0443928968868248
$ws.Range("F11").Value = 0.4331058859825134
$ws.Range("G11").Value = 0.7

# Row 12
$ws.Range("D12").Value = 0.0000001179622106730826
$ws.Range("E12").Value = 0.0000001179622106730826

# Row 13
$ws.Range("D13").Value = 0.9961414439456551
$ws.Range("E13").Value = 0.9961414439456551

# Row 14
$ws.Range("D14").Value = 0.0000004617682079374064
$ws.Range("E14").Value = 0.0000004617682079374064

# Row 15
$ws.Range("D15").Value = 0.000000000000000000000128430598019855
$ws.Range("E15").Value = 0.000000000000000000000128430598019855

# Row 16
$ws.Range("D16").Value = 0.5015821953256491
$ws.Range("E16").Value = 0.5015821953256491

# Row 17
$ws.Range("D17").Value = 0.9999999999999774
$ws.Range("E17").Value = 0.00000000000002264854970235319

# Row 18
$ws.Range("D18").Value = 0.9963149205326404
$ws.Range("E18").Value = 0.003685079467359631

# Row 19
$ws.Range("D19").Value = 0.9999999988078023
$ws.Range("E19").Value = 0.000000001192197673915985

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.08231867041990477
$ws.Range("E20").Value = 0.9176813295800952

# Row 21
$ws.Range("D21").Value = 0.9635051288125224
$ws.Range("E21").Value = 0.03649487118747763
$ws.Range("F21").Value = 0.8791804313659668
$ws.Range("G21").Value = 0.7

$wb.Save()
